$d = $word.ActiveDocument

function Set-ParaText($para, [string]$text) {
    $r = $para.Range
    $r.MoveEnd(1, -1) | Out-Null
    if ($r.Start -ne $r.End) { $r.Delete() }
    $ins = $para.Range
    $ins.MoveEnd(1, -1) | Out-Null
    $ins.Collapse(0)
    $ins.InsertAfter($text)
}

# Paragraph 1: title
Set-ParaText ($d.Paragraphs.Item(1)) 'ContosoLearn Market Research'

# Paragraph 2: AdatumLearn (merge 2 runs into 1)
Set-ParaText ($d.Paragraphs.Item(2)) 'AdatumLearn: AdatumLearn is a top AI-powered learning platform that uses artificial intelligence to enrich eLearning with features that automate a variety of tasks. It is known for its content authoring capabilities and adaptive learning technology.'

# Paragraph 3
Set-ParaText ($d.Paragraphs.Item(3)) 'AdventureLearn: AdventureLearn is another AI-powered learning platform that offers personalized learning experiences and data-driven recommendations.'

# Paragraph 4
Set-ParaText ($d.Paragraphs.Item(4)) 'AlpineTraining: AlpineTraining is a mobile-first learning platform that focuses on microlearning.'

# Paragraph 5
Set-ParaText ($d.Paragraphs.Item(5)) 'Bellows OnDemand: Bellows OnDemand is a comprehensive learning solution that offers content creation and social collaboration.'

# Paragraph 6
Set-ParaText ($d.Paragraphs.Item(6)) 'FabrikamLearning: FabrikamLearning provides a suite of learning platforms that cater to different learning needs.'

# Paragraph 7
Set-ParaText ($d.Paragraphs.Item(7)) 'FirstUp Cards: FirstUp Cards is a mobile learning app that is ideal for training on safety procedures, compliance, new product knowledge or any other type of training scenario.'

# Paragraph 8
Set-ParaText ($d.Paragraphs.Item(8)) 'Munson''sLearn: Munson''sLearn is designed to enable businesses to train their employees, partners, and customers.'

# Paragraph 9
Set-ParaText ($d.Paragraphs.Item(9)) 'LibertyLearn: LibertyLearn is a fast LMS for your mission-critical project.'

# Paragraph 10: WoodgroveLMS (split into 3 runs with proofErr markers)
$p10 = $d.Paragraphs.Item(10)
$r10 = $p10.Range
$r10.MoveEnd(1, -1) | Out-Null
if ($r10.Start -ne $r10.End) { $r10.Delete() }
$ins10 = $p10.Range
$ins10.MoveEnd(1, -1) | Out-Null
$ins10.Collapse(0)
$ins10.InsertAfter('WoodgroveLMS: WoodgroveLMS is a functional and attractive learning management system built to provide ')
$ins10.Collapse(0)
$ins10.InsertAfter('a best')
$ins10.Collapse(0)
$ins10.InsertAfter('-in-class training experience.')

# Paragraph 11
Set-ParaText ($d.Paragraphs.Item(11)) 'NorthwindWorlds: NorthwindWorlds is a powerful, easy-to-use, and reliable training solution for individuals and enterprises.'

# Paragraph 12
Set-ParaText ($d.Paragraphs.Item(12)) 'ProsewareLearn: ProsewareLearn is an online education company that offers a variety of video training courses for software developers, IT administrators, and creative professionals through its website.'

# Paragraph 13
Set-ParaText ($d.Paragraphs.Item(13)) 'RelecloudLearn: RelecloudLearn is an American online learning platform that offers massive open online courses (MOOC), specializations, and degrees in a variety of subjects.'

# Paragraph 14
Set-ParaText ($d.Paragraphs.Item(14)) 'TreyAcademy: TreyAcademy is an online learning platform aimed at professional adults and students, developed in May 2010.'

# Paragraph 15: merge 3 runs into 1 (final summary paragraph)
Set-ParaText ($d.Paragraphs.Item(15)) 'These platforms have a significant market presence and are widely recognized for their AI-powered features, such as personalized learning experiences, data-driven recommendations, and automation of tasks. They are transforming the eLearning landscape by leveraging AI to deliver more engaging, rewarding, and personalized learning experiences. '

Write-Host "Done"
